# Edit applies the "missions-counts-long" update:
# adds a second (regenerated) copy of the summary table below the first one
# (rows 23-42), mirroring the header row and populating refreshed
# counts/codes, consistent with a Power Query "data" refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the new header row (row 23) as a formatted copy of row 1 ---
$ws.Range("A1:O1").Copy()
$ws.Range("A23:O23").PasteSpecial(-4122)   # xlPasteFormats - carries over the header styles/borders

$ws.Cells.Item(23, 1).Value  = "1. Main Objective"
$ws.Cells.Item(23, 2).Value  = "num"
$ws.Cells.Item(23, 3).Value  = "codes"
$ws.Cells.Item(23, 4).Value  = "2. Side Quest"
$ws.Cells.Item(23, 5).Value  = "num"
$ws.Cells.Item(23, 6).Value  = "codes"
$ws.Cells.Item(23, 7).Value  = "3. Early Game"
$ws.Cells.Item(23, 8).Value  = "num"
$ws.Cells.Item(23, 9).Value  = "codes"
$ws.Cells.Item(23, 10).Value = "4. Mid Game"
$ws.Cells.Item(23, 11).Value = "num"
$ws.Cells.Item(23, 12).Value = "codes"
$ws.Cells.Item(23, 13).Value = "5. Late Game"
$ws.Cells.Item(23, 14).Value = "num"
$ws.Cells.Item(23, 15).Value = "codes"

# --- Populate the refreshed data rows (24-42) ---
$ws.Cells.Item(24, 1).Value = 'Open all 9 worlds'
$ws.Cells.Item(24, 2).Value = 8531
$ws.Cells.Item(24, 3).Value = 'J'
$ws.Cells.Item(24, 4).Value = 'All 3 Cheato Visits'
$ws.Cells.Item(24, 5).Value = 6879
$ws.Cells.Item(24, 7).Value = 'Begin run w/ MM 100% Trotless'
$ws.Cells.Item(24, 8).Value = 10966
$ws.Cells.Item(24, 9).Value = 'R'
$ws.Cells.Item(24, 10).Value = '13 tokens in RBB [r 10-15]'
$ws.Cells.Item(24, 11).Value = 6004
$ws.Cells.Item(24, 12).Value = 'TA'
$ws.Cells.Item(24, 13).Value = 'Both HCs in BGS'
$ws.Cells.Item(24, 14).Value = 6428
$ws.Cells.Item(24, 15).Value = 'HA'
$ws.Cells.Item(25, 1).Value = 'All notes'
$ws.Cells.Item(25, 2).Value = 8469
$ws.Cells.Item(25, 3).Value = 'NR'
$ws.Cells.Item(25, 4).Value = 'No FFM'
$ws.Cells.Item(25, 5).Value = 6872
$ws.Cells.Item(25, 6).Value = 'R'
$ws.Cells.Item(25, 7).Value = 'Termite''s Quest: 5 jiggies and 1 HC as the termite'
$ws.Cells.Item(25, 8).Value = 10946
$ws.Cells.Item(25, 10).Value = '80 notes in RBB [r 40-100]'
$ws.Cells.Item(25, 11).Value = 5969
$ws.Cells.Item(25, 12).Value = 'NA'
$ws.Cells.Item(25, 13).Value = 'All Jinjos in BGS'
$ws.Cells.Item(25, 14).Value = 6410
$ws.Cells.Item(25, 15).Value = 'OA'
$ws.Cells.Item(26, 1).Value = 'All 24 honeycombs'
$ws.Cells.Item(26, 2).Value = 8427
$ws.Cells.Item(26, 3).Value = 'HT'
$ws.Cells.Item(26, 4).Value = 'All 9 orange Jinjos (color randomly chosen)'
$ws.Cells.Item(26, 5).Value = 6853
$ws.Cells.Item(26, 6).Value = 'O'
$ws.Cells.Item(26, 7).Value = 'No more than 2 jiggies in MM'
$ws.Cells.Item(26, 8).Value = 10240
$ws.Cells.Item(26, 9).Value = 'JR'
$ws.Cells.Item(26, 10).Value = 'All Jinjos in RBB'
$ws.Cells.Item(26, 11).Value = 5960
$ws.Cells.Item(26, 12).Value = 'OA'
$ws.Cells.Item(26, 13).Value = 'Tiptup''s jiggy'
$ws.Cells.Item(26, 14).Value = 6358
$ws.Cells.Item(26, 15).Value = 'J'
$ws.Cells.Item(27, 1).Value = 'Humanitarian'
$ws.Cells.Item(27, 2).Value = 8384
$ws.Cells.Item(27, 3).Value = 'J'
$ws.Cells.Item(27, 4).Value = 'All 10 Brentilda visits'
$ws.Cells.Item(27, 5).Value = 6825
$ws.Cells.Item(27, 7).Value = 'All tokens in TTC'
$ws.Cells.Item(27, 8).Value = 7774
$ws.Cells.Item(27, 9).Value = 'TA'
$ws.Cells.Item(27, 10).Value = 'Both HCs in RBB'
$ws.Cells.Item(27, 11).Value = 5899
$ws.Cells.Item(27, 12).Value = 'HA'
$ws.Cells.Item(27, 13).Value = '88 notes in BGS [r 75-100]'
$ws.Cells.Item(27, 14).Value = 6318
$ws.Cells.Item(27, 15).Value = 'NA'
$ws.Cells.Item(28, 1).Value = 'Defeat Grunty'
$ws.Cells.Item(28, 2).Value = 8374
$ws.Cells.Item(28, 3).Value = 'N'
$ws.Cells.Item(28, 4).Value = 'Activate all 8 warp cauldrons (not Dingpot)'
$ws.Cells.Item(28, 5).Value = 6821
$ws.Cells.Item(28, 7).Value = 'Both HCs in TTC'
$ws.Cells.Item(28, 8).Value = 7707
$ws.Cells.Item(28, 9).Value = 'HA'
$ws.Cells.Item(28, 10).Value = '7 jiggies in RBB [r 4-10]'
$ws.Cells.Item(28, 11).Value = 5630
$ws.Cells.Item(28, 12).Value = 'JA'
$ws.Cells.Item(28, 13).Value = 'Both timed jiggies in BGS'
$ws.Cells.Item(28, 14).Value = 6291
$ws.Cells.Item(28, 15).Value = 'J'
$ws.Cells.Item(29, 1).Value = 'All 116 tokens'
$ws.Cells.Item(29, 2).Value = 8364
$ws.Cells.Item(29, 3).Value = 'T'
$ws.Cells.Item(29, 4).Value = 'No RBA'
$ws.Cells.Item(29, 5).Value = 6815
$ws.Cells.Item(29, 6).Value = 'R'
$ws.Cells.Item(29, 7).Value = '9 jiggies in TTC [r 8-10]'
$ws.Cells.Item(29, 8).Value = 7332
$ws.Cells.Item(29, 9).Value = 'JA'
$ws.Cells.Item(29, 10).Value = 'GV rings jiggy without flight or bee'
$ws.Cells.Item(29, 11).Value = 5592
$ws.Cells.Item(29, 13).Value = 'All tokens in BGS'
$ws.Cells.Item(29, 14).Value = 6225
$ws.Cells.Item(29, 15).Value = 'TA'
$ws.Cells.Item(30, 1).Value = 'Open All 12 Note Doors and Defeat Grunty'
$ws.Cells.Item(30, 2).Value = 8296
$ws.Cells.Item(30, 3).Value = 'NR'
$ws.Cells.Item(30, 4).Value = '87 tokens [r 70-90]'
$ws.Cells.Item(30, 5).Value = 6725
$ws.Cells.Item(30, 6).Value = 'T'
$ws.Cells.Item(30, 7).Value = 'Both HCs in CC'
$ws.Cells.Item(30, 8).Value = 4505
$ws.Cells.Item(30, 9).Value = 'HA'
$ws.Cells.Item(30, 10).Value = 'Both HCs in GV'
$ws.Cells.Item(30, 11).Value = 5315
$ws.Cells.Item(30, 12).Value = 'HA'
$ws.Cells.Item(30, 13).Value = '7 jiggies in BGS [r 5-8]'
$ws.Cells.Item(30, 14).Value = 6204
$ws.Cells.Item(30, 15).Value = 'J'
$ws.Cells.Item(31, 1).Value = 'All Jinjos'
$ws.Cells.Item(31, 2).Value = 8279
$ws.Cells.Item(31, 3).Value = 'O'
$ws.Cells.Item(31, 4).Value = 'Open the 640 note door'
$ws.Cells.Item(31, 5).Value = 6710
$ws.Cells.Item(31, 6).Value = 'N'
$ws.Cells.Item(31, 7).Value = 'All Jinjos in CC'
$ws.Cells.Item(31, 8).Value = 4434
$ws.Cells.Item(31, 9).Value = 'OA'
$ws.Cells.Item(31, 10).Value = 'Abuse Gobi (beak bust Gobi at all 5 locations)'
$ws.Cells.Item(31, 11).Value = 5293
$ws.Cells.Item(31, 12).Value = 'J'
$ws.Cells.Item(31, 13).Value = 'Croctuses jiggy'
$ws.Cells.Item(31, 14).Value = 6184
$ws.Cells.Item(31, 15).Value = 'J'
$ws.Cells.Item(32, 1).Value = 'Open 765 note door'
$ws.Cells.Item(32, 2).Value = 8253
$ws.Cells.Item(32, 3).Value = 'N'
$ws.Cells.Item(32, 4).Value = 'All lair jiggies'
$ws.Cells.Item(32, 5).Value = 6637
$ws.Cells.Item(32, 6).Value = 'J'
$ws.Cells.Item(32, 7).Value = 'All 4 jiggies inside Clanker'
$ws.Cells.Item(32, 8).Value = 4420
$ws.Cells.Item(32, 9).Value = 'J'
$ws.Cells.Item(32, 10).Value = 'All tokens in GV'
$ws.Cells.Item(32, 11).Value = 5246
$ws.Cells.Item(32, 12).Value = 'TA'
$ws.Cells.Item(32, 13).Value = 'All 21 caterpillars'
$ws.Cells.Item(32, 14).Value = 4734
$ws.Cells.Item(33, 1).Value = 'All of 1 type of collectible from each world'
$ws.Cells.Item(33, 2).Value = 8247
$ws.Cells.Item(33, 3).Value = 'A'
$ws.Cells.Item(33, 4).Value = 'All 5 transformations'
$ws.Cells.Item(33, 5).Value = 6631
$ws.Cells.Item(33, 6).Value = 'T'
$ws.Cells.Item(33, 7).Value = '89 notes in CC [r 80-100]'
$ws.Cells.Item(33, 8).Value = 4336
$ws.Cells.Item(33, 9).Value = 'NA'
$ws.Cells.Item(33, 10).Value = 'All Jinjos in GV'
$ws.Cells.Item(33, 11).Value = 5240
$ws.Cells.Item(33, 12).Value = 'OA'
$ws.Cells.Item(33, 13).Value = 'Kill all 6 Sir Slushes in winter'
$ws.Cells.Item(33, 14).Value = 4731
$ws.Cells.Item(34, 1).Value = '78 jiggies [r 75-90]'
$ws.Cells.Item(34, 2).Value = 8196
$ws.Cells.Item(34, 3).Value = 'JR'
$ws.Cells.Item(34, 4).Value = '2 jiggies from each world'
$ws.Cells.Item(34, 5).Value = 6542
$ws.Cells.Item(34, 6).Value = 'J'
$ws.Cells.Item(34, 7).Value = '10 jiggies in CC [r 8-10]'
$ws.Cells.Item(34, 8).Value = 4233
$ws.Cells.Item(34, 9).Value = 'JA'
$ws.Cells.Item(34, 10).Value = '87 notes in GV [r 40-100]'
$ws.Cells.Item(34, 11).Value = 5212
$ws.Cells.Item(34, 12).Value = 'NA'
$ws.Cells.Item(34, 13).Value = 'Eyrie''s jiggy'
$ws.Cells.Item(34, 14).Value = 4643
$ws.Cells.Item(35, 1).Value = 'Open DoG & defeat Grunty'
$ws.Cells.Item(35, 2).Value = 8180
$ws.Cells.Item(35, 3).Value = 'NJR'
$ws.Cells.Item(35, 4).Value = '40 jiggies [r 40-60]'
$ws.Cells.Item(35, 5).Value = 6514
$ws.Cells.Item(35, 6).Value = 'J'
$ws.Cells.Item(35, 7).Value = 'Both HCs in FP'
$ws.Cells.Item(35, 8).Value = 4020
$ws.Cells.Item(35, 9).Value = 'HA'
$ws.Cells.Item(35, 10).Value = 'Kill all 10 Limbos (skeletons) in MMM'
$ws.Cells.Item(35, 11).Value = 5120
$ws.Cells.Item(35, 13).Value = 'All Jinjos in CCW'
$ws.Cells.Item(35, 14).Value = 4591
$ws.Cells.Item(35, 15).Value = 'OA'
$ws.Cells.Item(36, 4).Value = 'No MMM early'
$ws.Cells.Item(36, 5).Value = 6475
$ws.Cells.Item(36, 6).Value = 'RJ'
$ws.Cells.Item(36, 7).Value = 'All Jinjos in FP'
$ws.Cells.Item(36, 8).Value = 3888
$ws.Cells.Item(36, 9).Value = 'OA'
$ws.Cells.Item(36, 10).Value = 'MMM witch switch jiggy'
$ws.Cells.Item(36, 11).Value = 5058
$ws.Cells.Item(36, 13).Value = 'Nabnut''s jiggy'
$ws.Cells.Item(36, 14).Value = 4523
$ws.Cells.Item(36, 15).Value = 'J'
$ws.Cells.Item(37, 4).Value = 'No FP early'
$ws.Cells.Item(37, 5).Value = 6447
$ws.Cells.Item(37, 6).Value = 'RJ'
$ws.Cells.Item(37, 7).Value = '99 notes in FP [r 80-100]'
$ws.Cells.Item(37, 8).Value = 3865
$ws.Cells.Item(37, 9).Value = 'NA'
$ws.Cells.Item(37, 10).Value = '8 jiggies in GV [r 4-9]'
$ws.Cells.Item(37, 11).Value = 4893
$ws.Cells.Item(37, 12).Value = 'JA'
$ws.Cells.Item(37, 13).Value = 'Flower jiggy in CCW'
$ws.Cells.Item(37, 14).Value = 4505
$ws.Cells.Item(37, 15).Value = 'J'
$ws.Cells.Item(38, 4).Value = '15 HCs [r 14-18]'
$ws.Cells.Item(38, 5).Value = 6254
$ws.Cells.Item(38, 6).Value = 'HA'
$ws.Cells.Item(38, 7).Value = 'Merry Christmas! (Visit Boggy''s igloo w/ him in it & give presents)'
$ws.Cells.Item(38, 8).Value = 3845
$ws.Cells.Item(38, 9).Value = 'J'
$ws.Cells.Item(38, 10).Value = 'Both HCs in MMM'
$ws.Cells.Item(38, 11).Value = 4878
$ws.Cells.Item(38, 12).Value = 'HA'
$ws.Cells.Item(38, 13).Value = '25 tokens in CCW [r 15-25]'
$ws.Cells.Item(38, 14).Value = 4480
$ws.Cells.Item(38, 15).Value = 'TA'
$ws.Cells.Item(39, 7).Value = 'All tokens in FP'
$ws.Cells.Item(39, 8).Value = 3831
$ws.Cells.Item(39, 9).Value = 'TA'
$ws.Cells.Item(39, 10).Value = 'All Jinjos in MMM'
$ws.Cells.Item(39, 11).Value = 4801
$ws.Cells.Item(39, 12).Value = 'OA'
$ws.Cells.Item(39, 13).Value = 'Both HCs in CCW'
$ws.Cells.Item(39, 14).Value = 4467
$ws.Cells.Item(39, 15).Value = 'HA'
$ws.Cells.Item(40, 7).Value = '9 jiggies in FP [r 4-9]'
$ws.Cells.Item(40, 8).Value = 3658
$ws.Cells.Item(40, 9).Value = 'JA'
$ws.Cells.Item(40, 10).Value = '68 notes in MMM [r 60-100]'
$ws.Cells.Item(40, 11).Value = 4747
$ws.Cells.Item(40, 12).Value = 'NA'
$ws.Cells.Item(40, 13).Value = '69 notes in CCW [r 50-80]'
$ws.Cells.Item(40, 14).Value = 4434
$ws.Cells.Item(40, 15).Value = 'NA'
$ws.Cells.Item(41, 10).Value = '14 tokens in MMM [r 10-16]'
$ws.Cells.Item(41, 11).Value = 4698
$ws.Cells.Item(41, 12).Value = 'TA'
$ws.Cells.Item(41, 13).Value = 'Collect 8 jiggies as the bee'
$ws.Cells.Item(41, 14).Value = 4330
$ws.Cells.Item(41, 15).Value = 'JTR'
$ws.Cells.Item(42, 10).Value = '6 jiggies in MMM [r 6-10]'
$ws.Cells.Item(42, 11).Value = 4445
$ws.Cells.Item(42, 12).Value = 'JA'
$ws.Cells.Item(42, 13).Value = '6 jiggies in CCW [r 4-8]'
$ws.Cells.Item(42, 14).Value = 4144
$ws.Cells.Item(42, 15).Value = 'JA'

# --- Re-apply the descending sort on the new "5. Late Game" mini-table
#     (M24:O42 sorted by N24:N42 descending), matching the sortState Excel
#     records after a Power Query refresh + re-sort ---
$sortRange = $ws.Range("M24:O42")
$sortKey = $ws.Range("N24:N42")
$sortRange.Sort($sortKey, 2)

# --- Update the active selection/scroll position to match where the user
#     ended up after adding the table ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("Q32").Select()
